$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string text values used in column B
$ws.Range("B7").Value = "First Entry"
$ws.Range("B8").Value = "Second Entry"

# Update the selected range/active cell in the sheet view
$ws.Range("B9").Select()
